$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 171047
$ws.Range("C4").Value = 161859
$ws.Range("C5").Value = 9189
$ws.Range("C8").Value = 65.88
